$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.788608551025391
$ws.Range("B1").Value = 4.140255928039551
$ws.Range("C1").Value = 2.023570775985718
$ws.Range("D1").Value = 0.8889397382736206
$ws.Range("E1").Value = 0.479514479637146
